# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 9244
$ws1.Range("F7").Value = 11974
$ws1.Range("F8").Value = 11974
$ws1.Range("F22").Value = 311
$ws1.Range("F35").Value = 502
$ws1.Range("F37").Value = 534
$ws1.Range("F39").Value = 2183
$ws1.Range("F44").Value = 445
$ws1.Range("F46").Value = 885
$ws1.Range("F50").Value = 270

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F25").Value = 426

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 261

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 9244
$ws4.Range("F10").Value = 11974
$ws4.Range("F20").Value = 311
$ws4.Range("F26").Value = 261
$ws4.Range("F35").Value = 502
$ws4.Range("F37").Value = 534
$ws4.Range("F39").Value = 2183
$ws4.Range("F43").Value = 445
$ws4.Range("F47").Value = 426
$ws4.Range("F50").Value = 270
